$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 0.484733
$ws.Range("N2").Value = 1.454199
$ws.Range("O2").Value = 0.00792098608860474
$ws.Range("P2").Value = 0.00792098608860474
$ws.Range("Q2").Value = 55.85160002105467
$ws.Range("R2").Value = 502.664400189492
$ws.Range("S2").Value = 0.002207902588879627
$ws.Range("T2").Value = 0.002207902588879626

$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.1147190689515559
$ws.Range("P3").Value = 0.1147190689515559
$ws.Range("Q3").Value = 808.8946858633747
$ws.Range("R3").Value = 7280.052172770372
$ws.Range("S3").Value = 0.03197689359616294
$ws.Range("T3").Value = 0.03197689359616294

$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("M4").Value = 53.289524
$ws.Range("N4").Value = 159.868572
$ws.Range("O4").Value = 0.8708001689019901
$ws.Range("P4").Value = 0.8708001689019901
$ws.Range("Q4").Value = 6140.09192640153
$ws.Range("R4").Value = 55260.82733761377
$ws.Range("S4").Value = 0.2427276005548683
$ws.Range("T4").Value = 0.2427276005548683

$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 0.4014323333333333
$ws.Range("N5").Value = 1.204297
$ws.Range("O5").Value = 0.006559776057849319
$ws.Range("P5").Value = 0.006559776057849319
$ws.Range("Q5").Value = 46.25358314134178
$ws.Range("R5").Value = 416.282248272076
$ws.Range("S5").Value = 0.001828477714590622
$ws.Range("T5").Value = 0.001828477714590621

$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 0.484733
$ws.Range("N6").Value = 1.454199
$ws.Range("O6").Value = 0.00792098608860474
$ws.Range("P6").Value = 0.00792098608860474
$ws.Range("Q6").Value = 89.598921208866
$ws.Range("R6").Value = 806.390290879794
$ws.Range("S6").Value = 0.003541987875428843
$ws.Range("T6").Value = 0.003541987875428842

$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.1147190689515559
$ws.Range("P7").Value = 0.1147190689515559
$ws.Range("S7").Value = 0.05129835436669362
$ws.Range("T7").Value = 0.05129835436669362

$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("M8").Value = 53.289524
$ws.Range("N8").Value = 159.868572
$ws.Range("O8").Value = 0.8708001689019901
$ws.Range("P8").Value = 0.8708001689019901
$ws.Range("Q8").Value = 9850.131643882249
$ws.Range("R8").Value = 88651.18479494023
$ws.Range("S8").Value = 0.3893913719416139
$ws.Range("T8").Value = 0.3893913719416138

$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 0.4014323333333333
$ws.Range("N9").Value = 1.204297
$ws.Range("O9").Value = 0.006559776057849319
$ws.Range("P9").Value = 0.006559776057849319
$ws.Range("Q9").Value = 74.201475874398
$ws.Range("R9").Value = 667.8132828695819
$ws.Range("S9").Value = 0.002933302369493672
$ws.Range("T9").Value = 0.002933302369493672

$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 0.484733
$ws.Range("N10").Value = 1.454199
$ws.Range("O10").Value = 0.00792098608860474
$ws.Range("P10").Value = 0.00792098608860474
$ws.Range("Q10").Value = 29.35104687790334
$ws.Range("R10").Value = 264.15942190113
$ws.Range("S10").Value = 0.001160293570168455
$ws.Range("T10").Value = 0.001160293570168455

$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.1147190689515559
$ws.Range("P11").Value = 0.1147190689515559
$ws.Range("Q11").Value = 425.0890902877034
$ws.Range("R11").Value = 3825.80181258933
$ws.Range("S11").Value = 0.01680444790474925
$ws.Range("T11").Value = 0.01680444790474925

$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("M12").Value = 53.289524
$ws.Range("N12").Value = 159.868572
$ws.Range("O12").Value = 0.8708001689019901
$ws.Range("P12").Value = 0.8708001689019901
$ws.Range("Q12").Value = 3226.731658511294
$ws.Range("R12").Value = 29040.58492660164
$ws.Range("S12").Value = 0.1275578350443183
$ws.Range("T12").Value = 0.1275578350443183

$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 0.4014323333333333
$ws.Range("N13").Value = 1.204297
$ws.Range("O13").Value = 0.006559776057849319
$ws.Range("P13").Value = 0.006559776057849319
$ws.Range("Q13").Value = 24.30711182026556
$ws.Range("R13").Value = 218.76400638239
$ws.Range("S13").Value = 0.0009608987942318482
$ws.Range("T13").Value = 0.0009608987942318481

$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 0.484733
$ws.Range("N14").Value = 1.454199
$ws.Range("O14").Value = 0.00792098608860474
$ws.Range("P14").Value = 0.00792098608860474
$ws.Range("Q14").Value = 25.56947589624167
$ws.Range("R14").Value = 230.125283066175
$ws.Range("S14").Value = 0.001010802054127815
$ws.Range("T14").Value = 0.001010802054127815

$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.1147190689515559
$ws.Range("P15").Value = 0.1147190689515559
$ws.Range("Q15").Value = 370.3208711117417
$ws.Range("R15").Value = 3332.887840005675
$ws.Range("S15").Value = 0.01463937308395007
$ws.Range("T15").Value = 0.01463937308395007

$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("M16").Value = 53.289524
$ws.Range("N16").Value = 159.868572
$ws.Range("O16").Value = 0.8708001689019901
$ws.Range("P16").Value = 0.8708001689019901
$ws.Range("Q16").Value = 2811.001519269767
$ws.Range("R16").Value = 25299.0136734279
$ws.Range("S16").Value = 0.1111233613611896
$ws.Range("T16").Value = 0.1111233613611896

$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 0.4014323333333333
$ws.Range("N17").Value = 1.204297
$ws.Range("O17").Value = 0.006559776057849319
$ws.Range("P17").Value = 0.006559776057849319
$ws.Range("Q17").Value = 21.17539835566944
$ws.Range("R17").Value = 190.578585201025
$ws.Range("S17").Value = 0.0008370971795331763
$ws.Range("T17").Value = 0.0008370971795331763
